# edit.ps1 - applies the "update with new config file" edit to workflow.pptx
#
# Summary of the change (per the target diff):
#  1. The big background rectangle on slide 2 ("CasellaDiTesto 3") switches
#     its fill from accent4 to accent2 and loses its tx1 outline (noFill),
#     while keeping the 28575 EMU line weight.
#  2. The small label box on slide 2 ("CasellaDiTesto 5") that used to read
#     "MIABIS compliant dataset" on one line now reads "INPUT" / "Dataset
#     XLSX" on two centered lines.
#  3. The cached datetimeFigureOut field text (footer date) is refreshed
#     from 18/04/2024 to 04/07/2024 across the slide master and every
#     slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1 & 2: shape restyle + text rewrite on slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape 1 = "CasellaDiTesto 3", the big rounded rectangle behind everything.
$bigRect = $s2.Shapes.Item(1)
$bigRect.Fill.ForeColor.ObjectThemeColor = 6   # msoThemeColorAccent2
$bigRect.Line.Visible = 0                      # drop the tx1 outline -> noFill

# Shape 3 = "CasellaDiTesto 5", the small callout that names the dataset.
$labelBox = $s2.Shapes.Item(3)
$labelBox.TextFrame.TextRange.Text = "INPUT" + [char]13 + "Dataset XLSX"

# ---------------------------------------------------------------------
# 3: refresh the cached footer date everywhere it is placeholders-driven
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "04/07/2024"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes "04/07/2024"
}
